# SCD0204 - Validasi Mockup Digisales Mobile and
# SCD0205 - Validasi Field report PHR pada searching portal
#
# The PHR report screen gained two new search/result fields ("TEXT4" /
# "TEXT5") right before the existing "FILE1" column. We make room for them
# by moving the old column O ("FILE1" header / file-name value) two columns
# to the right (to Q), filling the vacated O/P header cells with the two new
# field names, and narrowing O/P to the same width used by the other
# short "TEXTn" columns while widening the relocated column back to its
# original width.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the source widths before anything moves, so the relocated column
# keeps its original (wide) width and the vacated columns can match the
# other short "TEXTn" columns.
$narrowWidth = $ws.Columns("L").ColumnWidth
$wideWidth = $ws.Columns("O").ColumnWidth

# --- Move the "FILE1" column (header + sample value) from O to Q -----------
# Copy preserves both value and number/style formatting of the source cell.
$ws.Range("O1").Copy($ws.Range("Q1"))
$ws.Range("O2").Copy($ws.Range("Q2"))

# Clear the old column O cells completely (no leftover cell/style at O2,
# matching the fully-vacated source column).
$ws.Range("O1").ClearContents()
$ws.Range("O2").Clear()

# --- Fill in the two new header labels --------------------------------------
$ws.Range("O1").Value2 = "TEXT4"
$ws.Range("P1").Value2 = "TEXT5"

# --- Column widths ------------------------------------------------------
# O/P become narrow "TEXTn"-style columns (same width as L:M), Q takes on
# the width the old column O used to have.
$ws.Columns("O:P").ColumnWidth = $narrowWidth
$ws.Columns("Q").ColumnWidth = $wideWidth

# --- View / selection -----------------------------------------------------
# Scroll the sheet so column G is the left-most visible column and select
# the relocated "FILE1" value cell.
$excel.ActiveWindow.ScrollColumn = 7
$ws.Range("Q2").Select() | Out-Null
